$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-99 down to 10-100.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new data record.
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44750
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112052
$ws.Range("G9").Value = "Albahaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 3500
$ws.Range("L9").Value = 4000
$ws.Range("M9").Value = 3750
$ws.Range("N9").Value = "$/paquete"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 3750
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
